$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.970.03"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.58%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.581.54"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.59%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.03"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.67"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.97%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.59"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.51%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.76%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.350"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.05"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.66%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.045.46"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.880.40"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.56%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.579.80"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.11"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "340.56"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.34"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.91%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.64"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.31%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.43"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.01%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +7.41%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.53%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.01"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.26"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.94"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.77%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "460.03"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0800"

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.56%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "176.35"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.25%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.397"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.87"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.48"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.09%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.70"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "158.52"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.43%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.69%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.69"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.41%  "

$ws.Range("B44").Value = "InjectiveProtocol"

$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.30"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.45%  "

$ws.Range("B45").Value = "Mantle"

$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.636"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0538"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0961"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0236"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.98%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.97"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.40"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.18%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.45%  "
